# "Actualizar" refresh: shifts the history of availability-check timestamps
# down one slot and stamps the newest refresh time into the top block.
#
# Column D holds a custom-formatted (YYYY-MM-DD HH:MM:SS) date/time serial
# number. Each block of 14 rows shares the same timestamp, and this update
# replaces:
#   D2:D15  44241.48606363745 -> 44241.50730187182  (new "now" timestamp)
#   D16:D29 44241.46486099537 -> 44241.48606363426  (previous block's value)
#   D30:D43 44241.44365131945 -> 44241.46486099537  (previous block's value)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value  = 44241.50730187182
$ws.Range("D16:D29").Value = 44241.48606363426
$ws.Range("D30:D43").Value = 44241.46486099537
